$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New column C values (set in this order so the new shared-string
#     entries line up with the target indices 32-38) ---
$ws.Range("C3").Value  = "914-IN-U11-U011-1"
$ws.Range("C5").Value  = "D1"
$ws.Range("C6").Value  = "Station04"
$ws.Range("C7").Value  = "04v1"
$ws.Range("C8").Value  = "04v2"
$ws.Range("C9").Value  = "Station05"
$ws.Range("C10").Value = "05v1"

# --- New column D values for rows 9 & 10 (shared-string indices 39-40) ---
$ws.Range("D9").Value  = "Bh1st2"
$ws.Range("D10").Value = "st2v1"

# --- New font formatting applied only to C6 (Source Sans Pro, #333333) ---
$ws.Range("C6").Font.Color = 3355443
$ws.Range("C6").Font.Name = "Source Sans Pro"

# --- Move the active selection from C18 to C17 ---
$ws.Range("C17").Select()

# --- Window size bookkeeping recorded by Excel on save ---
$excel.ActiveWindow.Width = 23085
$excel.ActiveWindow.Height = 10845
